$wb = $excel.ActiveWorkbook

# --- Sheet1: update the active selection (was B7, now B9) ---
$sheet1 = $wb.Worksheets.Item("Sheet1")
[void]$sheet1.Range("B9").Select()

# --- Add new Sheet2 right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Name = "Sheet2"

# Column widths to roughly match the authored layout
$ws2.Columns.Item(1).ColumnWidth = 30.25
$ws2.Columns.Item(2).ColumnWidth = 14.75
$ws2.Columns.Item(3).ColumnWidth = 15.75
$ws2.Columns.Item(4).ColumnWidth = 13.25

# --- Populate cells in the exact order the values were authored, so that
#     new shared-string entries come out in the same sequence as the source
#     workbook (Username, Password, the new login, then the remaining
#     address/header fields). ---
$ws2.Range("A1").Value = "Username"
$ws2.Range("B1").Value = "Password"

$ws2.Range("A3").Value = "hitendravibhandik1@gmail.com"
$ws2.Range("B3").Value = "Hitesh@222"

$ws2.Range("C1").Value = "Street Address"
$ws2.Range("D1").Value = "Apt no"

$ws2.Range("C3").Value = "carmel drive"
$ws2.Range("D2").Value = "A453"
$ws2.Range("D3").Value = "B305"

$ws2.Range("E1").Value = "Guests"
$ws2.Range("C2").Value = "rochester hills"

# --- Row 2 (existing credentials reused from Sheet1) ---
$ws2.Range("A2").Value = "vijeyata.thorat@gmail.com"
$ws2.Range("B2").Value = "Vijeyata@12345"

# --- Guest counts ---
$ws2.Range("E2").Value = 3
$ws2.Range("E3").Value = 4

$ws2.Hyperlinks.Add($ws2.Range("A3"), "mailto:hitendravibhandik1@gmail.com")
$ws2.Hyperlinks.Add($ws2.Range("B3"), "mailto:Hitesh@222")

# --- Make Sheet2 the active sheet/tab with its own selection ---
$ws2.Activate()
[void]$ws2.Range("C6").Select()
